# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker table (rows 16-33, cols C/D/E/F) is re-sorted/rewritten:
#  - the old row 16 (JOSE GREGORIO ALCAZAR ARRIETA / 9101444 / period 2005)
#    is dropped from the top and re-appended at the bottom (row 33).
#  - SENEN MARTINEZ TORREGLOSA's periods are re-ordered newest-first
#    (2110 down to 2005) occupying rows 16-32.
#  - a brand new period (2110) is added for SENEN at row 16, carrying the
#    older "Salario Basico" start date (29260) that used to sit on the last
#    row; the rest of SENEN's rows keep the 35112 start date.
#  - JOSE's re-appended row at the bottom reverts to the 35112 start date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docId = "8853287"
$name = "SENEN MARTINEZ TORREGLOSA"

$periods = @("2110","2109","2108","2107","2106","2105","2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2005")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 3).Value = $docId
    $ws.Cells.Item($row, 4).Value = $name
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($row -eq 16) {
        $ws.Cells.Item($row, 6).Value = 29260
    } else {
        $ws.Cells.Item($row, 6).Value = 35112
    }
}

# Row 33 reverts to the original top entry (JOSE GREGORIO ALCAZAR ARRIETA)
$ws.Cells.Item(33, 3).Value = "9101444"
$ws.Cells.Item(33, 4).Value = "JOSE GREGORIO ALCAZAR ARRIETA"
$ws.Cells.Item(33, 5).Value = "2005"
$ws.Cells.Item(33, 6).Value = 35112
